$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 ("Pianificazione errata") - Riscontro effettivo: "Nessuno" -> "Sì"
$ws.Range("C9").Value = "Sì"

# Update the visible selection to the full table range
$ws.Range("A1:C11").Select()
